# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect newly scraped totals (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listing)
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    4  = 721
    5  = 52
    7  = 2805
    8  = 1675
    9  = 1809
    12 = 734
    13 = 890
    14 = 166
    15 = 371
    16 = 1113
    18 = 49
    20 = 6562
    21 = 252
    22 = 1495
    23 = 156
    26 = 307
    27 = 263
    28 = 64
    30 = 901
    32 = 89
    34 = 473
    35 = 1425
    37 = 140
    39 = 12
    41 = 193
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Sheet "全部类型" (all-types combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    4  = 721
    5  = 52
    10 = 2805
    11 = 1675
    12 = 1809
    15 = 734
    17 = 890
    18 = 166
    19 = 371
    20 = 1113
    21 = 49
    23 = 6562
    24 = 252
    25 = 1496
    27 = 156
    30 = 307
    31 = 263
    32 = 64
    34 = 901
    36 = 89
    38 = 473
    39 = 1425
    41 = 140
    43 = 12
    45 = 193
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
